{"js": "// Append a trailing space (same bold / 15pt formatting as the heading\n// text) right after the \"Project Description\" heading, matching the\n// author's edit recorded in the diff:\n//   <w:t>Project Description</w:t> ... + <w:t xml:space=\"preserve\"> </w:t>\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Project Description\" heading paragraph (exact text match,\n// so we don't accidentally touch \"Project Purpose\"/\"Project Objectives\"\n// or any other heading).\nlet heading = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Project Description\") {\n    heading = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!heading) {\n  throw new Error('Could not find the \"Project Description\" paragraph.');\n}\n\n// Insert the extra space at the very end of that paragraph's text. The\n// new text inherits the run formatting already present at that\n// location (bold, 15pt/size 30 half-points), exactly like the run added\n// in the diff.\nheading.getRange(\"End\").insertText(\" \", \"End\");\n\nawait context.sync();\n", "ps1": "# Append a trailing space (inheriting the bold / 15pt heading\n# formatting) right after the \"Project Description\" heading text,\n# matching the author's edit:\n#   <w:t>Project Description</w:t> ... + <w:t xml:space=\"preserve\"> </w:t>\n\n$d = $word.ActiveDocument\n\n# Locate the exact \"Project Description\" heading text in the document.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Project Description\")\n\nif (-not $found) {\n    throw 'Could not find the \"Project Description\" text.'\n}\n\n# $rng now spans exactly \"Project Description\". Collapse to its end\n# point and insert the extra space right there; the new text picks up\n# the formatting already in effect at that point (bold, size 15 =\n# w:sz 30), same as the run added in the diff.\n$rng.Collapse(0)\n$rng.InsertAfter(\" \")\n"}
